$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column BF ("Date") was stamped with the export file name
# ("5-20-2013-14") instead of the actual game date. NBA.com's stats
# feed reports a game night under the following calendar day, so the
# correct ISO date for this export is 2014-05-20. Fix rows 2-31.
#
# Force Text format first so Excel keeps the literal "2014-05-20"
# string instead of silently re-parsing it back into a date serial.
$dateRange = $ws.Range("BF2:BF31")
$dateRange.NumberFormat = "@"

for ($row = 2; $row -le 31; $row++) {
    $ws.Cells.Item($row, 58).Value = "2014-05-20"
}
